$wb = $excel.ActiveWorkbook

# Data for new row 55 on each sheet: B,C,D,E are text; F,G,H,I are numeric; A is a date-time numeric with style copied from A54
$rowsData = @{
    "FE_LFT_#1" = @{ A=45841.49278935185; B="0x01,0x7c"; C="0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D="0x01,0x54"; E="0xf";  F=380; G=[double]"7.598631275147109e+23"; H=340; I=15 }
    "FE_LFT_#2" = @{ A=45841.49278935185; B="0x01,0x90"; C="0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D="0x01,0x64"; E="0xe";  F=400; G=[double]"5.68432987514711e+23";  H=356; I=14 }
    "FE_PLT_#1" = @{ A=45841.49278935185; B="0x00,0x6e"; C="0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D="0x00,0x68"; E="0x3";  F=110; G=[double]"5.68631262647114e+23";  H=104; I=3 }
    "FE_PLT_#2" = @{ A=45841.49278935185; B="0x00,0x6e"; C="0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D="0x00,0x68"; E="0x3";  F=110; G=[double]"9.85046333984776e+23";  H=104; I=3 }
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rowsData[$sheetName]
    $newRow = 55

    # Column A: numeric date value, copy the style/number format from the row above
    $ws.Cells.Item($newRow, 1).Value2 = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    # Columns B-E: inline text strings
    $ws.Cells.Item($newRow, 2).Value2 = $data.B
    $ws.Cells.Item($newRow, 3).Value2 = $data.C
    $ws.Cells.Item($newRow, 4).Value2 = $data.D
    $ws.Cells.Item($newRow, 5).Value2 = $data.E

    # Columns F-I: numeric values
    $ws.Cells.Item($newRow, 6).Value2 = $data.F
    $ws.Cells.Item($newRow, 7).Value2 = $data.G
    $ws.Cells.Item($newRow, 8).Value2 = $data.H
    $ws.Cells.Item($newRow, 9).Value2 = $data.I
}
